$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 4076.5715
$ws_ALC.Range("J17").Value = 4076.5715
$ws_ALC.Range("L17").Value = 12229.7145
$ws_ALC.Range("N17").Value = -12565.7145

# ALC row 43
$ws_ALC.Range("H43").Value = 1663.5555
$ws_ALC.Range("I43").Value = 1667.5
$ws_ALC.Range("K43").Value = 1667.5
$ws_ALC.Range("M43").Value = -1598.5

# ALC row 74
$ws_ALC.Range("H74").Value = 7343.05
$ws_ALC.Range("I74").Value = 4098.5
$ws_ALC.Range("J74").Value = 7703.5557
$ws_ALC.Range("K74").Value = 4098.5
$ws_ALC.Range("L74").Value = 7703.5557
$ws_ALC.Range("M74").Value = -3162.5
$ws_ALC.Range("N74").Value = -9575.555700000001

# ALC row 77
$ws_ALC.Range("H77").Value = 7343.05
$ws_ALC.Range("I77").Value = 4098.5
$ws_ALC.Range("J77").Value = 7703.5557
$ws_ALC.Range("K77").Value = 20492.5
$ws_ALC.Range("L77").Value = 38517.7785
$ws_ALC.Range("M77").Value = -15812.5
$ws_ALC.Range("N77").Value = -47877.7785

# ALC row 98
$ws_ALC.Range("H98").Value = 1392.6
$ws_ALC.Range("I98").Value = 1233.7241
$ws_ALC.Range("K98").Value = 1233.7241
$ws_ALC.Range("M98").Value = 264.2759000000001

# ALC row 113
$ws_ALC.Range("H113").Value = 6184
$ws_ALC.Range("J113").Value = 6299.75
$ws_ALC.Range("L113").Value = 6299.75
$ws_ALC.Range("N113").Value = -12807.75

# ALC row 121
$ws_ALC.Range("H121").Value = 2228.5293
$ws_ALC.Range("J121").Value = 2228.5293
$ws_ALC.Range("L121").Value = 6685.5879
$ws_ALC.Range("N121").Value = -10179.5879

# ALC row 122
$ws_ALC.Range("H122").Value = 1392.6
$ws_ALC.Range("I122").Value = 1233.7241
$ws_ALC.Range("K122").Value = 3701.1723
$ws_ALC.Range("M122").Value = -1251.1723

# ALC row 125
$ws_ALC.Range("H125").Value = 2997.1785
$ws_ALC.Range("I125").Value = 1666.4286
$ws_ALC.Range("J125").Value = 3440.762
$ws_ALC.Range("K125").Value = 14997.8574
$ws_ALC.Range("L125").Value = 30966.858
$ws_ALC.Range("M125").Value = -12537.8574
$ws_ALC.Range("N125").Value = -35886.858

# ALC row 131
$ws_ALC.Range("H131").Value = 6168.8
$ws_ALC.Range("I131").Value = 4665.1665
$ws_ALC.Range("J131").Value = 8424.25
$ws_ALC.Range("K131").Value = 13995.4995
$ws_ALC.Range("L131").Value = 25272.75
$ws_ALC.Range("M131").Value = -8955.499500000002
$ws_ALC.Range("N131").Value = -35352.75

# ALC row 137
$ws_ALC.Range("H137").Value = 2855
$ws_ALC.Range("I137").Value = 2412.0625
$ws_ALC.Range("J137").Value = 3248.7222
$ws_ALC.Range("K137").Value = 7236.1875
$ws_ALC.Range("L137").Value = 9746.1666
$ws_ALC.Range("M137").Value = -4686.1875
$ws_ALC.Range("N137").Value = -14846.1666

# ARM row 32
$ws_ARM.Range("H32").Value = 3596.8384
$ws_ARM.Range("I32").Value = 2229.926
$ws_ARM.Range("J32").Value = 9747.944
$ws_ARM.Range("K32").Value = 2229.926
$ws_ARM.Range("L32").Value = 9747.944
$ws_ARM.Range("M32").Value = -1942.926
$ws_ARM.Range("N32").Value = -10321.944

# ARM row 34
$ws_ARM.Range("H34").Value = 0
$ws_ARM.Range("I34").Value = 0
$ws_ARM.Range("K34").Value = 0
$ws_ARM.Range("M34").ClearContents()

# ARM row 42
$ws_ARM.Range("H42").Value = 10000
$ws_ARM.Range("J42").Value = 10000
$ws_ARM.Range("L42").Value = 10000
$ws_ARM.Range("N42").Value = -10972

# ARM row 97
$ws_ARM.Range("H97").Value = 32213.223
$ws_ARM.Range("I97").Value = 36101.8
$ws_ARM.Range("J97").Value = 27352.5
$ws_ARM.Range("K97").Value = 36101.8
$ws_ARM.Range("L97").Value = 27352.5
$ws_ARM.Range("M97").Value = -35605.8
$ws_ARM.Range("N97").Value = -28344.5

# ARM row 102
$ws_ARM.Range("H102").Value = 5112.567
$ws_ARM.Range("I102").Value = 4707.154
$ws_ARM.Range("K102").Value = 4707.154
$ws_ARM.Range("M102").Value = -3085.154

# ARM row 132
$ws_ARM.Range("H132").Value = 2037.8125
$ws_ARM.Range("I132").Value = 1559.3448
$ws_ARM.Range("K132").Value = 4678.0344
$ws_ARM.Range("M132").Value = -2148.0344

# CRP row 5
$ws_CRP.Range("H5").Value = 2558.3333
$ws_CRP.Range("I5").Value = 5000
$ws_CRP.Range("J5").Value = 2070
$ws_CRP.Range("K5").Value = 5000
$ws_CRP.Range("L5").Value = 2070
$ws_CRP.Range("M5").Value = -4888
$ws_CRP.Range("N5").Value = -2294

# CRP row 31
$ws_CRP.Range("H31").Value = 3046.1462
$ws_CRP.Range("I31").Value = 1883.3928
$ws_CRP.Range("K31").Value = 1883.3928
$ws_CRP.Range("M31").Value = -1588.3928

# CRP row 34
$ws_CRP.Range("H34").Value = 3046.1462
$ws_CRP.Range("I34").Value = 1883.3928
$ws_CRP.Range("K34").Value = 1883.3928
$ws_CRP.Range("M34").Value = -1681.3928

# CRP row 99
$ws_CRP.Range("H99").Value = 3991.75
$ws_CRP.Range("I99").Value = 3352.0908
$ws_CRP.Range("J99").Value = 5399
$ws_CRP.Range("K99").Value = 3352.0908
$ws_CRP.Range("L99").Value = 5399
$ws_CRP.Range("M99").Value = -1854.0908
$ws_CRP.Range("N99").Value = -8395

# CRP row 122
$ws_CRP.Range("H122").Value = 3601.9443
$ws_CRP.Range("J122").Value = 4250.75
$ws_CRP.Range("L122").Value = 12752.25
$ws_CRP.Range("N122").Value = -17652.25

# CRP row 126
$ws_CRP.Range("H126").Value = 3991.75
$ws_CRP.Range("I126").Value = 3352.0908
$ws_CRP.Range("J126").Value = 5399
$ws_CRP.Range("K126").Value = 10056.2724
$ws_CRP.Range("L126").Value = 16197
$ws_CRP.Range("M126").Value = -7586.2724
$ws_CRP.Range("N126").Value = -21137

# CRP row 132
$ws_CRP.Range("H132").Value = 57803.055
$ws_CRP.Range("I132").Value = 57803.055
$ws_CRP.Range("J132").Value = 0
$ws_CRP.Range("K132").Value = 173409.165
$ws_CRP.Range("L132").Value = 0
$ws_CRP.Range("M132").Value = -170879.165
$ws_CRP.Range("N132").ClearContents()

# CUL row 11
$ws_CUL.Range("H11").Value = 14999.857
$ws_CUL.Range("I11").Value = 20000
$ws_CUL.Range("K11").Value = 60000
$ws_CUL.Range("M11").Value = -59860

# CUL row 113
$ws_CUL.Range("H113").Value = 3574.8462
$ws_CUL.Range("I113").Value = 5710
$ws_CUL.Range("K113").Value = 17130
$ws_CUL.Range("M113").Value = -14960

# GSM row 2
$ws_GSM.Range("H2").Value = 55.625
$ws_GSM.Range("I2").Value = 57.5
$ws_GSM.Range("J2").Value = 53.75
$ws_GSM.Range("K2").Value = 57.5
$ws_GSM.Range("L2").Value = 53.75
$ws_GSM.Range("M2").Value = 55.5
$ws_GSM.Range("N2").Value = -279.75

# GSM row 11
$ws_GSM.Range("H11").Value = 0
$ws_GSM.Range("J11").Value = 0
$ws_GSM.Range("L11").Value = 0
$ws_GSM.Range("N11").ClearContents()

# GSM row 62
$ws_GSM.Range("H62").Value = 514999.5
$ws_GSM.Range("I62").Value = 514999.5
$ws_GSM.Range("J62").Value = 0
$ws_GSM.Range("K62").Value = 514999.5
$ws_GSM.Range("L62").Value = 0
$ws_GSM.Range("M62").Value = -514313.5
$ws_GSM.Range("N62").ClearContents()

# GSM row 65
$ws_GSM.Range("H65").Value = 514999.5
$ws_GSM.Range("I65").Value = 514999.5
$ws_GSM.Range("J65").Value = 0
$ws_GSM.Range("K65").Value = 1544998.5
$ws_GSM.Range("L65").Value = 0
$ws_GSM.Range("M65").Value = -1541566.5
$ws_GSM.Range("N65").ClearContents()

# GSM row 70
$ws_GSM.Range("H70").Value = 6283.3335
$ws_GSM.Range("I70").Value = 6283.3335
$ws_GSM.Range("K70").Value = 6283.3335
$ws_GSM.Range("M70").Value = -6013.3335

# GSM row 73
$ws_GSM.Range("H73").Value = 6283.3335
$ws_GSM.Range("I73").Value = 6283.3335
$ws_GSM.Range("K73").Value = 6283.3335
$ws_GSM.Range("M73").Value = -5347.3335

# GSM row 102
$ws_GSM.Range("H102").Value = 81098.62
$ws_GSM.Range("I102").Value = 3646.3333
$ws_GSM.Range("J102").Value = 104334.3
$ws_GSM.Range("K102").Value = 3646.3333
$ws_GSM.Range("L102").Value = 104334.3
$ws_GSM.Range("M102").Value = -2024.3333
$ws_GSM.Range("N102").Value = -107578.3

# GSM row 122
$ws_GSM.Range("H122").Value = 70577.11
$ws_GSM.Range("I122").Value = 104365.5
$ws_GSM.Range("K122").Value = 313096.5
$ws_GSM.Range("M122").Value = -310646.5

# GSM row 132
$ws_GSM.Range("H132").Value = 6999
$ws_GSM.Range("I132").Value = 6999
$ws_GSM.Range("K132").Value = 20997
$ws_GSM.Range("M132").Value = -18467

# LTW row 12
$ws_LTW.Range("H12").Value = 0
$ws_LTW.Range("J12").Value = 0
$ws_LTW.Range("L12").Value = 0
$ws_LTW.Range("N12").ClearContents()

# LTW row 38
$ws_LTW.Range("H38").Value = 18000
$ws_LTW.Range("J38").Value = 0
$ws_LTW.Range("L38").Value = 0
$ws_LTW.Range("N38").ClearContents()

# LTW row 68
$ws_LTW.Range("H68").Value = 3499.75
$ws_LTW.Range("I68").Value = 3999.5
$ws_LTW.Range("K68").Value = 3999.5
$ws_LTW.Range("M68").Value = -3250.5

# LTW row 71
$ws_LTW.Range("H71").Value = 3499.75
$ws_LTW.Range("I71").Value = 3999.5
$ws_LTW.Range("K71").Value = 19997.5
$ws_LTW.Range("M71").Value = -16253.5

# LTW row 96
$ws_LTW.Range("H96").Value = 30992
$ws_LTW.Range("J96").Value = 30992
$ws_LTW.Range("L96").Value = 30992
$ws_LTW.Range("N96").Value = -36484

# LTW row 136
$ws_LTW.Range("H136").Value = 60801.484
$ws_LTW.Range("J136").Value = 7181.636
$ws_LTW.Range("L136").Value = 21544.908
$ws_LTW.Range("N136").Value = -26644.908

# WVR row 10
$ws_WVR.Range("H10").Value = 0
$ws_WVR.Range("I10").Value = 0
$ws_WVR.Range("K10").Value = 0
$ws_WVR.Range("M10").ClearContents()

# WVR row 41
$ws_WVR.Range("H41").Value = 77499.57000000001
$ws_WVR.Range("J41").Value = 77499.57000000001
$ws_WVR.Range("L41").Value = 77499.57000000001
$ws_WVR.Range("N41").Value = -78279.57000000001

# WVR row 96
$ws_WVR.Range("H96").Value = 1865.6666
$ws_WVR.Range("J96").Value = 1865.6666
$ws_WVR.Range("L96").Value = 1865.6666
$ws_WVR.Range("N96").Value = -4611.6666

# WVR row 122
$ws_WVR.Range("H122").Value = 3295.9443
$ws_WVR.Range("I122").Value = 2497.0715
$ws_WVR.Range("J122").Value = 6092
$ws_WVR.Range("K122").Value = 7491.2145
$ws_WVR.Range("L122").Value = 18276
$ws_WVR.Range("M122").Value = -5041.2145
$ws_WVR.Range("N122").Value = -23176

# WVR row 126
$ws_WVR.Range("H126").Value = 2556.95
$ws_WVR.Range("I126").Value = 2617.4285
$ws_WVR.Range("J126").Value = 2415.8333
$ws_WVR.Range("K126").Value = 7852.2855
$ws_WVR.Range("L126").Value = 7247.499899999999
$ws_WVR.Range("M126").Value = -5382.2855
$ws_WVR.Range("N126").Value = -12187.4999

# WVR row 132
$ws_WVR.Range("H132").Value = 319667.62
$ws_WVR.Range("I132").Value = 8025.174
$ws_WVR.Range("J132").Value = 1116087.2
$ws_WVR.Range("K132").Value = 24075.522
$ws_WVR.Range("L132").Value = 3348261.6
$ws_WVR.Range("M132").Value = -21545.522
$ws_WVR.Range("N132").Value = -3353321.6
